# Feria Lagunitas de Puerto Montt - Poroto verde: weekly data refresh.
# A new observation is inserted as row 44 (pushing the previously-existing
# rows 44..133 down to 45..134); dimension grows from A1:R133 to A1:R134.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 44, shifting everything below it down.
$ws.Rows.Item(44).Insert()

# Populate the new row 44 with the latest weekly record.
$ws.Cells.Item(44, 1).Value  = 4
$ws.Cells.Item(44, 2).Value  = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(44, 3).Value  = 'Los Lagos'
$ws.Cells.Item(44, 4).Value  = 45002
$ws.Cells.Item(44, 5).Value  = 10
$ws.Cells.Item(44, 6).Value  = 100112031
$ws.Cells.Item(44, 7).Value  = 'Poroto verde'
$ws.Cells.Item(44, 8).Value  = 'Magnum'
$ws.Cells.Item(44, 9).Value  = 'Primera'
$ws.Cells.Item(44, 10).Value = 40
$ws.Cells.Item(44, 11).Value = 33000
$ws.Cells.Item(44, 12).Value = 33000
$ws.Cells.Item(44, 13).Value = 33000
$ws.Cells.Item(44, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(44, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(44, 16).Value = 1320
$ws.Cells.Item(44, 17).Value = 25
$ws.Cells.Item(44, 18).Value = 'Hortaliza'
